$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets

$s = $ws.Item(1)
$s.Range("F2").Value = 2343
$s.Range("F3").Value = 509
$s.Range("F5").Value = 342
$s.Range("F6").Value = 342
$s.Range("F7").Value = 557
$s.Range("F9").Value = 770
$s.Range("F10").Value = 536
$s.Range("F11").Value = 792
$s.Range("F12").Value = 383
$s.Range("F13").Value = 96
$s.Range("F14").Value = 393
$s.Range("F15").Value = 19
$s.Range("F16").Value = 1024
$s.Range("F17").Value = 20551
$s.Range("F18").Value = 708
$s.Range("F19").Value = 70
$s.Range("F20").Value = 242
$s.Range("F21").Value = 292
$s.Range("F22").Value = 174
$s.Range("F23").Value = 150
$s.Range("F26").Value = 221
$s.Range("F27").Value = 19
$s.Range("F28").Value = 339
$s.Range("F29").Value = 150

$s = $ws.Item(2)
$s.Range("F2").Value = 35
$s.Range("F3").Value = 183
$s.Range("G3").Value = 480
$s.Range("F4").Value = 2
$s.Range("F5").Value = 88
$s.Range("F7").Value = 223
$s.Range("F8").Value = 3413
$s.Range("F10").Value = 94
$s.Range("F14").Value = 123
$s.Range("F16").Value = 3762
$s.Range("G16").Value = "已售罄"

$s = $ws.Item(3)
$s.Range("F2").Value = 267
$s.Range("F3").Value = 99
$s.Range("F4").Value = 594
$s.Range("F5").Value = 206

$s = $ws.Item(4)
$s.Range("F2").Value = 267
$s.Range("F3").Value = 99
$s.Range("F4").Value = 35
$s.Range("F5").Value = 2343
$s.Range("F6").Value = 594
$s.Range("F7").Value = 509
$s.Range("F9").Value = 342
$s.Range("F10").Value = 342
$s.Range("F11").Value = 557
$s.Range("F12").Value = 183
$s.Range("G12").Value = 480
$s.Range("F14").Value = 2
$s.Range("F15").Value = 88
$s.Range("F17").Value = 206
$s.Range("F18").Value = 770
$s.Range("F19").Value = 536
$s.Range("F20").Value = 792
$s.Range("F21").Value = 383
$s.Range("F22").Value = 96
$s.Range("F23").Value = 393
$s.Range("F25").Value = 1024
$s.Range("F26").Value = 20552
$s.Range("F27").Value = 223
$s.Range("F28").Value = 3413
$s.Range("F30").Value = 94
$s.Range("F32").Value = 708
$s.Range("F33").Value = 70
$s.Range("F34").Value = 242
$s.Range("F37").Value = 292
$s.Range("F38").Value = 174
$s.Range("F39").Value = 150
$s.Range("F42").Value = 123
$s.Range("F44").Value = 221
$s.Range("F45").Value = 19
$s.Range("F46").Value = 339
$s.Range("F47").Value = 150
$s.Range("F48").Value = 3762
$s.Range("G48").Value = "已售罄"